$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.363.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.655.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.70"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.260"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("E10").Value = "  -1.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0876"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.888.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.652.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.83%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.570"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.66"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.365.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0729"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.25%  "

$ws.Range("E20").Value = "  -0.53%  "

$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.98%  "

$ws.Range("E27").Value = "  -1.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("E31").Value = "  -3.98%  "

$ws.Range("E32").Value = "  -1.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.457.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.18%  "

$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("E35").Value = "  -0.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("E37").Value = "  -2.40%  "

$ws.Range("E38").Value = "  -1.23%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("E40").Value = "  +0.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("E42").Value = "  +0.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.29%  "

$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.797.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.63%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.786"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.52%  "

$ws.Range("E47").Value = "  +0.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("E49").Value = "  -3.27%  "

$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.42%  "
